# The weekly update adds one new price observation for the current week.
# It is inserted as a new row 57 (pushing all subsequent rows down by one),
# matching the existing table's row layout/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57 (row 1 is the header, rows 2..152 are data).
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(57, 1).Value  = 8
$ws.Cells.Item(57, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(57, 3).Value  = "Coquimbo"
$ws.Cells.Item(57, 4).Value  = 44477
$ws.Cells.Item(57, 5).Value  = 4
$ws.Cells.Item(57, 6).Value  = 100112003
$ws.Cells.Item(57, 7).Value  = "Ajo"
$ws.Cells.Item(57, 8).Value  = "Chino"
$ws.Cells.Item(57, 9).Value  = "Primera"
$ws.Cells.Item(57, 10).Value = 600
$ws.Cells.Item(57, 11).Value = 15000
$ws.Cells.Item(57, 12).Value = 16000
$ws.Cells.Item(57, 13).Value = 15500
$ws.Cells.Item(57, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(57, 15).Value = "China"
$ws.Cells.Item(57, 16).Value = 1550
$ws.Cells.Item(57, 17).Value = 10
$ws.Cells.Item(57, 18).Value = "Hortaliza"

# Give the date cell the same number format (style index 2) used by the
# rest of the "Fecha" column.
$ws.Cells.Item(57, 4).NumberFormat = $ws.Cells.Item(58, 4).NumberFormat
